# Update specific numeric values in the KNN imputation result sheet
# (commit message: "Update Name of Algo")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C3"  = -12.708
    "B7"  = 5.525
    "A10" = -21.252
    "D10" = -7.926
    "A12" = -21.515
    "D14" = -7.637
    "B15" = 5.243
    "A18" = -21.866
    "C18" = -11.081
    "C19" = -12.051
    "B20" = 7.217000000000001
    "C27" = -13.022
    "B29" = 5.626
    "B30" = 5.149
    "B31" = 6.114000000000001
    "D32" = -8.274000000000001
    "D35" = -7.387
    "A37" = -20.02
    "B40" = 9.208
    "C42" = -12.328
    "D43" = -8.375
    "C44" = -12.324
    "C47" = -12.357
    "D49" = -7.997000000000002
    "A55" = -21.875
    "D56" = -8.087
    "C58" = -12.674
    "A68" = -21.581
    "B68" = 5.242
    "D69" = -6.812
    "C73" = -12.77
    "B76" = 6.654000000000001
    "A77" = -20.252
    "A78" = -19.533
    "D81" = -7.342000000000001
    "B87" = 5.429
    "B88" = 5.75
    "D92" = -6.587999999999999
    "C95" = -11.663
    "B96" = 6.484
    "B98" = 5.140999999999999
    "B101" = 8.659000000000001
    "C101" = -12.643
    "B102" = 7.739999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
